$wb = $excel.ActiveWorkbook

# --- currency_conversions sheet: split "foreign_amount" into explicit
# source_amount / target_amount, and add a matching target_fees column ---
$ws = $wb.Worksheets.Item("currency_conversions")

# Insert two new columns after source_currency (D) for target_amount and
# target_fees, pushing the old target_currency / comment columns right
$ws.Range("E1:F1").EntireColumn.Insert()

# New header row: target_fees, then rename "foreign_amount" to
# "source_amount", then add target_amount (order chosen to match the
# resulting shared-string table layout)
$ws.Range("F1").Value = "target_fees"
$ws.Range("B1").Value = "source_amount"
$ws.Range("E1").Value = "target_amount"

# Fill the new target_amount / target_fees values for the existing rows
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = 0
$ws.Range("E3").Value = -1
$ws.Range("F3").Value = 0

# Make this the active sheet/tab
$ws.Activate()
